$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark (it currently sits in the
# "The probability of failure ..." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}

# Delete the duplicated block of paragraphs (second copy of the
# "Assignment4.4" writeup through the trailing blank paragraphs),
# keeping only the very last (empty) paragraph before the section break.
$startPara = $d.Paragraphs.Item(37)
$endPara = $d.Paragraphs.Item(111)
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()

# Re-add the "_GoBack" bookmark into the final (now last) empty paragraph.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
